$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    12262,10951,10554,10554,10554,10554,10336,10336,10336,10336,
    10336,9949,9949,9949,9549,8760,8760,8760,8760,8208,
    8028,8028,8028,8028,8028,8028,8028,8028,8028,8028,
    8028,8028,8028,8028,7903,7903,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

